# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh updates to the Tonberry_Profits workbook
# (columns H-N on the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 180.28572
$ws.Range("I55").Value = 202.5
$ws.Range("J55").Value = 150.66667
$ws.Range("K55").Value = 202.5
$ws.Range("L55").Value = 150.66667
$ws.Range("M55").Value = 11.5
$ws.Range("N55").Value = -578.6666700000001

$ws.Range("H74").Value = 4099.6
$ws.Range("I74").Value = 2749.5
$ws.Range("K74").Value = 2749.5
$ws.Range("M74").Value = -1813.5

$ws.Range("H77").Value = 4099.6
$ws.Range("I77").Value = 2749.5
$ws.Range("K77").Value = 13747.5
$ws.Range("M77").Value = -9067.5

$ws.Range("H116").Value = 20665.834
$ws.Range("J116").Value = 4799
$ws.Range("L116").Value = 4799
$ws.Range("N116").Value = -11683

$ws.Range("H135").Value = 233.6875
$ws.Range("I135").Value = 113.53333
$ws.Range("K135").Value = 1021.79997
$ws.Range("M135").Value = 1513.20003

$ws.Range("H138").Value = 1849.8727
$ws.Range("J138").Value = 2294.2334
$ws.Range("L138").Value = 6882.7002
$ws.Range("N138").Value = -17162.7002

$ws.Range("H141").Value = 3113142.2
$ws.Range("I141").Value = 3501660
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 10504980
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = -10499800
$ws.Range("N141").Value = -25360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 293566.53
$ws.Range("J2").Value = 1280.8
$ws.Range("L2").Value = 1280.8
$ws.Range("N2").Value = -1506.8

$ws.Range("H32").Value = 5647.6274
$ws.Range("I32").Value = 3653.6592
$ws.Range("J32").Value = 18181.143
$ws.Range("K32").Value = 3653.6592
$ws.Range("L32").Value = 18181.143
$ws.Range("M32").Value = -3366.6592
$ws.Range("N32").Value = -18755.143

$ws.Range("H34").Value = 40025
$ws.Range("I34").Value = 40025
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 40025
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -39754
$ws.Range("N34").ClearContents()

$ws.Range("H61").Value = 5515.0356
$ws.Range("I61").Value = 6424.294
$ws.Range("K61").Value = 6424.294
$ws.Range("M61").Value = -6212.294

$ws.Range("H74").Value = 1095.742
$ws.Range("I74").Value = 463.5238
$ws.Range("K74").Value = 463.5238
$ws.Range("M74").Value = 410.4762

$ws.Range("H77").Value = 1095.742
$ws.Range("I77").Value = 463.5238
$ws.Range("K77").Value = 2317.619
$ws.Range("M77").Value = 2050.381

$ws.Range("H116").Value = 293566.53
$ws.Range("J116").Value = 1280.8
$ws.Range("L116").Value = 1280.8
$ws.Range("N116").Value = -5868.8

$ws.Range("H132").Value = 1651.4193
$ws.Range("I132").Value = 1278.4584
$ws.Range("K132").Value = 3835.3752
$ws.Range("M132").Value = -1305.3752

$ws.Range("H136").Value = 5515.0356
$ws.Range("I136").Value = 6424.294
$ws.Range("K136").Value = 19272.882
$ws.Range("M136").Value = -16722.882

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 293566.53
$ws.Range("J3").Value = 1280.8
$ws.Range("L3").Value = 1280.8
$ws.Range("N3").Value = -1508.8

$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()

$ws.Range("H99").Value = 1588
$ws.Range("I99").Value = 1178.5
$ws.Range("J99").Value = 1997.5
$ws.Range("K99").Value = 1178.5
$ws.Range("L99").Value = 1997.5
$ws.Range("M99").Value = 319.5
$ws.Range("N99").Value = -4993.5

$ws.Range("H134").Value = 5542.2593
$ws.Range("I134").Value = 5892.087
$ws.Range("K134").Value = 17676.261
$ws.Range("M134").Value = -15141.261

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1822.2222
$ws.Range("I22").Value = 1000
$ws.Range("K22").Value = 1000
$ws.Range("M22").Value = -650

$ws.Range("H58").Value = 3107390.2
$ws.Range("J58").Value = 3008.8572
$ws.Range("L58").Value = 3008.8572
$ws.Range("N58").Value = -3414.8572

$ws.Range("H80").Value = 23148
$ws.Range("J80").Value = 27500
$ws.Range("L80").Value = 27500
$ws.Range("N80").Value = -29746

$ws.Range("H83").Value = 23148
$ws.Range("J83").Value = 27500
$ws.Range("L83").Value = 82500
$ws.Range("N83").Value = -93732

$ws.Range("H103").Value = 18144
$ws.Range("I103").Value = 17333
$ws.Range("K103").Value = 17333
$ws.Range("M103").Value = -16161

$ws.Range("H107").Value = 882.9375
$ws.Range("I107").Value = 652.8
$ws.Range("J107").Value = 1266.5
$ws.Range("K107").Value = 652.8
$ws.Range("L107").Value = 1266.5
$ws.Range("M107").Value = 1267.2
$ws.Range("N107").Value = -5106.5

$ws.Range("H132").Value = 2555.0952
$ws.Range("I132").Value = 1546.3572
$ws.Range("K132").Value = 4639.071599999999
$ws.Range("M132").Value = -2109.071599999999

$ws.Range("H134").Value = 3023.1428
$ws.Range("I134").Value = 2669.8
$ws.Range("K134").Value = 8009.400000000001
$ws.Range("M134").Value = -5474.400000000001

$ws.Range("H136").Value = 3107390.2
$ws.Range("J136").Value = 3008.8572
$ws.Range("L136").Value = 9026.571599999999
$ws.Range("N136").Value = -14126.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 455.45456
$ws.Range("J107").Value = 451
$ws.Range("L107").Value = 1353
$ws.Range("N107").Value = -5193

$ws.Range("H116").Value = 2331.3333
$ws.Range("I116").Value = 1037.25
$ws.Range("J116").Value = 2978.375
$ws.Range("K116").Value = 3111.75
$ws.Range("L116").Value = 8935.125
$ws.Range("M116").Value = 330.25
$ws.Range("N116").Value = -15819.125

$ws.Range("H131").Value = 10159.459
$ws.Range("J131").Value = 10508.464
$ws.Range("L131").Value = 31525.392
$ws.Range("N131").Value = -41605.392

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 23363.637
$ws.Range("I46").Value = 19500
$ws.Range("K46").Value = 19500
$ws.Range("M46").Value = -19344

$ws.Range("H70").Value = 3999.6667
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 3999.6667
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 3999.6667
$ws.Range("N70").Value = -4539.6667
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 3999.6667
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 3999.6667
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 3999.6667
$ws.Range("N73").Value = -5871.6667
$ws.Range("M73").ClearContents()

$ws.Range("H102").Value = 3063.5625
$ws.Range("I102").Value = 2875.5454
$ws.Range("K102").Value = 2875.5454
$ws.Range("M102").Value = -1253.5454

$ws.Range("H122").Value = 1483.1111
$ws.Range("I122").Value = 1306.8462
$ws.Range("K122").Value = 3920.5386
$ws.Range("M122").Value = -1470.5386

$ws.Range("H126").Value = 3144738.8
$ws.Range("J126").Value = 113543.11
$ws.Range("L126").Value = 340629.33
$ws.Range("N126").Value = -345569.33

$ws.Range("H132").Value = 1674389.4
$ws.Range("I132").Value = 2263833.8
$ws.Range("K132").Value = 6791501.399999999
$ws.Range("M132").Value = -6788971.399999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H96").Value = 75000
$ws.Range("J96").Value = 75000
$ws.Range("L96").Value = 75000
$ws.Range("N96").Value = -80492

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 15000
$ws.Range("J21").Value = 15000
$ws.Range("L21").Value = 15000
$ws.Range("N21").Value = -15470

$ws.Range("H35").Value = 15000
$ws.Range("J35").Value = 15000
$ws.Range("L35").Value = 15000
$ws.Range("N35").Value = -15580

$ws.Range("H122").Value = 28590.143
$ws.Range("I122").Value = 32988.793
$ws.Range("K122").Value = 98966.37899999999
$ws.Range("M122").Value = -96516.37899999999

$ws.Range("H126").Value = 2135.4546
$ws.Range("I126").Value = 2228
$ws.Range("K126").Value = 6684
$ws.Range("M126").Value = -4214

$ws.Range("H136").Value = 24156118
$ws.Range("I136").Value = 42736296
$ws.Range("J136").Value = 1889.7
$ws.Range("K136").Value = 128208888
$ws.Range("L136").Value = 5669.1
$ws.Range("M136").Value = -128206338
$ws.Range("N136").Value = -10769.1

